$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the M column values (rows 4-8) from 1.16 to 1.2
$ws.Range("M4").Value = 1.2
$ws.Range("M5").Value = 1.2
$ws.Range("M6").Value = 1.2
$ws.Range("M7").Value = 1.2
$ws.Range("M8").Value = 1.2

# Update the selected cell/range shown in the sheet view
$ws.Range("N12").Select()
